$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level changes -------------------------------------------------

# Rename the sheet tab.
$ws.Name = "Sheet2"

# Drop the stray _FilterDatabase defined names (workbook.xml <definedNames>).
$nameCount = $wb.Names.Count
for ($i = 1; $i -le $nameCount; $i++) {
    $wb.Names.Item(1).Delete()
}

# --- Worksheet cleanup -------------------------------------------------------

# The sheet carried 8 stray, empty "AMJ" cells (one per data row) that were
# inflating the used range all the way out to column AMJ (1024). Clearing
# them lets the sheet's real dimension collapse back down to A1:AW9.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1024).ClearContents()
}

# Re-flatten the column widths: instead of several distinct per-range widths,
# every column (1 through 1025, matching the old max) now shares one width.
$firstCol = $ws.Cells.Item(1, 1)
$lastCol = $ws.Cells.Item(1, 1025)
$ws.Range($firstCol, $lastCol).EntireColumn.ColumnWidth = 8.3

# Zoom to 100% and move the selection to F13.
$excel.ActiveWindow.Zoom = 100
$ws.Range("F13").Select()

# --- Page setup ---------------------------------------------------------

$ps = $ws.PageSetup
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7
$ps.CenterHeader = "&""Times New Roman,Regular""&12&A"
$ps.CenterFooter = "&""Times New Roman,Regular""&12Page &P"

Write-Host "done"
